$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 178930
$ws.Range("C4").Value = 168887
$ws.Range("C7").Value = 5.61
$ws.Range("C8").Value = 64.93000000000001
